$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '91.441.00'
$ws.Range("E2").Value = '  -3.72%  '
$ws.Range("D2:E2").ClearFormats()

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '3.304.88'
$ws.Range("E3").Value = '  -5.01%  '
$ws.Range("D3:E3").ClearFormats()

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '228.81'
$ws.Range("E5").Value = '  -4.50%  '
$ws.Range("D5:E5").ClearFormats()

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = '610.16'
$ws.Range("E6").Value = '  -5.34%  '
$ws.Range("D6:E6").ClearFormats()

# Row 7
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = '1.38'
$ws.Range("E7").Value = '  -5.35%  '
$ws.Range("D7:E7").ClearFormats()

# Row 8
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = '0.380'
$ws.Range("E8").Value = '  -6.26%  '
$ws.Range("D8:E8").ClearFormats()

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = '0.938'
$ws.Range("E10").Value = '  -6.30%  '
$ws.Range("D10:E10").ClearFormats()

# Row 11
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = '3.305.03'
$ws.Range("E11").Value = '  -4.94%  '
$ws.Range("D11:E11").ClearFormats()

# Row 12
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = '41.55'
$ws.Range("E12").Value = '  -2.72%  '
$ws.Range("D12:E12").ClearFormats()

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = '5.94'
$ws.Range("E14").Value = '  -4.15%  '
$ws.Range("D14:E14").ClearFormats()

# Row 15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = '91.317.39'
$ws.Range("E15").Value = '  -3.61%  '
$ws.Range("D15:E15").ClearFormats()

# Row 16
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = '3.928.56'
$ws.Range("E16").Value = '  -4.98%  '
$ws.Range("D16:E16").ClearFormats()

# Row 17
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000241'
$ws.Range("E17").Value = '  -5.94%  '
$ws.Range("D17:E17").ClearFormats()

# Row 18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = '8.01'
$ws.Range("E18").Value = '  -5.84%  '
$ws.Range("D18:E18").ClearFormats()

# Row 19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = '3.309.86'
$ws.Range("E19").Value = '  -5.28%  '
$ws.Range("D19:E19").ClearFormats()

# Row 20
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = '17.10'
$ws.Range("E20").Value = '  -5.14%  '
$ws.Range("D20:E20").ClearFormats()

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = '10.73'
$ws.Range("E21").Value = '  -6.67%  '
$ws.Range("D21:E21").ClearFormats()

# Row 22
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = '3.39'
$ws.Range("E22").Value = '  +6.09%  '
$ws.Range("D22:E22").ClearFormats()

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = '486.59'
$ws.Range("E23").Value = '  -3.78%  '
$ws.Range("D23:E23").ClearFormats()

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = '0.441'
$ws.Range("E24").Value = '  -14.49%  '
$ws.Range("D24:E24").ClearFormats()

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000178'
$ws.Range("E25").Value = '  -7.86%  '
$ws.Range("D25:E25").ClearFormats()

# Row 26
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = '6.14'
$ws.Range("E26").Value = '  -7.84%  '
$ws.Range("D26:E26").ClearFormats()

# Row 27
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = '88.81'
$ws.Range("E27").Value = '  -6.86%  '
$ws.Range("D27:E27").ClearFormats()

# Row 28
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = '11.68'
$ws.Range("E28").Value = '  -3.85%  '
$ws.Range("D28:E28").ClearFormats()

# Row 29
$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = '3.486.78'
$ws.Range("E29").Value = '  -4.96%  '
$ws.Range("D29:E29").ClearFormats()

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("B31:E31").NumberFormat = "@"
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '0.138'
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("B31:E31").ClearFormats()

# Row 32
$ws.Range("B32:E32").NumberFormat = "@"
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '10.86'
$ws.Range("E32").Value = '  -8.86%  '
$ws.Range("B32:E32").ClearFormats()

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -6.26%  '
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D34:E34").ClearFormats()

# Row 35
$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("D35").Value = '0.170'
$ws.Range("E35").Value = '  -7.44%  '
$ws.Range("D35:E35").ClearFormats()

# Row 36
$ws.Range("D36:E36").NumberFormat = "@"
$ws.Range("D36").Value = '27.70'
$ws.Range("E36").Value = '  -10.64%  '
$ws.Range("D36:E36").ClearFormats()

# Row 37
$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = '0.517'
$ws.Range("E37").Value = '  -9.84%  '
$ws.Range("D37:E37").ClearFormats()

# Row 38
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = '533.91'
$ws.Range("E38").Value = '  -2.14%  '
$ws.Range("D38:E38").ClearFormats()

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '7.27'
$ws.Range("E40").Value = '  -6.58%  '
$ws.Range("D40:E40").ClearFormats()

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.46%  '
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = '1.35'
$ws.Range("E42").Value = '  -7.88%  '
$ws.Range("D42:E42").ClearFormats()

# Row 43
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = '0.852'
$ws.Range("E43").Value = '  -9.60%  '
$ws.Range("D43:E43").ClearFormats()

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.43%  '
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '3.59'
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("D45:E45").ClearFormats()

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.02%  '
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0404'
$ws.Range("E47").Value = '  -3.05%  '
$ws.Range("D47:E47").ClearFormats()

# Row 48
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = '5.30'
$ws.Range("E48").Value = '  -6.88%  '
$ws.Range("D48:E48").ClearFormats()

# Row 49
$ws.Range("B49:E49").NumberFormat = "@"
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '2.06'
$ws.Range("E49").Value = '  -4.82%  '
$ws.Range("B49:E49").ClearFormats()

# Row 50
$ws.Range("B50:E50").NumberFormat = "@"
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = '51.43'
$ws.Range("E50").Value = '  -3.68%  '
$ws.Range("B50:E50").ClearFormats()

# Row 51
$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = '7.90'
$ws.Range("E51").Value = '  -2.27%  '
$ws.Range("D51:E51").ClearFormats()
